{"js": "// --- 1. Split \"Version\" into \"Versi\" | \"on\" --------------------------------\n// Dropping a zero-width bookmark at the boundary and immediately deleting it\n// forces Word to break the run at that character offset without leaving any\n// residual run formatting, matching the run split produced by the original\n// edit history.\nconst versiResults = context.document.body.search(\"Versi\", { matchCase: true });\nversiResults.load(\"items\");\nawait context.sync();\n\nconst afterVersi = versiResults.items[0].getRange(\"After\");\nafterVersi.insertBookmark(\"__TmpSplit\");\nawait context.sync();\n\ncontext.document.deleteBookmark(\"__TmpSplit\");\nawait context.sync();\n\n// --- 2. Bump the version number: \"1.\" -> \"2\" --------------------------------\n// The trailing period is dropped here; it is reinserted after the _GoBack\n// bookmark in step 3, so it ends up on the far side of the bookmark - exactly\n// like the target revision.\nconst numResults = context.document.body.search(\"1.\", { matchCase: true });\nnumResults.load(\"items\");\nawait context.sync();\n\nnumResults.items[0].insertText(\"2\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 3. Re-append the \".\" after the _GoBack bookmark -----------------------\nconst goBackRange = context.document.getBookmarkRange(\"_GoBack\");\ngoBackRange.insertText(\".\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Split \"Version\" into \"Versi\" | \"on\" --------------------------------\n# Dropping a zero-width bookmark at the boundary and immediately deleting it\n# forces Word to break the run at that character offset without leaving any\n# residual run formatting (<w:rPr/>), matching the run split produced by the\n# original edit history.\n$verRange = $d.Content\n$verRange.Find.ClearFormatting()\n$verRange.Find.MatchCase = $true\n$verRange.Find.Execute(\"Version\") | Out-Null\n$splitPos = $verRange.Start + 5\n$d.Bookmarks.Add(\"__TmpSplit\", $d.Range($splitPos, $splitPos)) | Out-Null\n$d.Bookmarks(\"__TmpSplit\").Delete()\n\n# --- 2. Bump the version number: \"1.\" -> \"2\" --------------------------------\n# The trailing period is dropped here; it is reinserted after the _GoBack\n# bookmark in step 3, so it ends up on the far side of the bookmark - exactly\n# like the target revision.\n$numRange = $d.Content\n$numRange.Find.ClearFormatting()\n$numRange.Find.MatchCase = $true\n$numRange.Find.Execute(\"1.\") | Out-Null\n$numRange.Text = \"2\"\n\n# --- 3. Re-append the \".\" after the _GoBack bookmark -----------------------\n$goBack = $d.Bookmarks(\"_GoBack\")\n$goBack.Range.InsertAfter(\".\")\n"}
